# IST price update 2025-12-22 14:48
#
# A new price-check timestamp column is inserted right after the "SKU Name"
# column (i.e. before the current column B), shifting every existing date
# column (B:Z) one slot to the right (-> C:AA). The new column B is filled
# with the latest scraped price for each SKU (falling back to the next
# known price when a SKU had no price recorded in the column that is now
# shifting into C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column B; this shifts B:Z -> C:AA,
# carries the dimension/col-width metadata along, and extends the used
# range to A1:AA26 automatically.
$ws.Columns.Item(2).Insert()

# The insert resets column B to the sheet default width; match it back up
# with the (identical) width the rest of the date columns use.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# New header timestamp for the freshly inserted column.
$ws.Range("B1").Value = "2025-12-22 20:14"

# Latest price snapshot for each SKU row (carry the prior latest price
# forward when a row had no price recorded yet).
$prices = @{
  2  = 929
  3  = 569
  4  = 299
  5  = 569
  6  = 499
  7  = 569
  8  = 929
  9  = 299
  10 = 299
  11 = 2997
  12 = 569
  13 = 569
  14 = 794
  15 = 499
  16 = 299
  17 = 929
  18 = 499
  19 = 1299
  20 = 929
  21 = 499
  22 = 299
  23 = 1299
  24 = 929
  25 = 929
  26 = 1299
}

foreach ($row in $prices.Keys) {
    $ws.Cells.Item($row, 2).Value = $prices[$row]
}
